$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D='''24.383.91'; E='  -2.12%  '},
    @{Row=3; D='''1.647.13'; E='  -3.82%  '},
    @{Row=4; D='''1.004'; E='  +0.05%  '},
    @{Row=5; D='''310.32'; E='  -0.96%  '},
    @{Row=6; D='''1.000'; E='  +0.32%  '},
    @{Row=7; D='''0.3650'; E='  -3.11%  '},
    @{Row=8; D='''46.55'; E='  -6.45%  '},
    @{Row=9; D='''0.3228'; E='  -7.03%  '},
    @{Row=10; D='''1.116'; E='  -8.47%  '},
    @{Row=11; D='''0.06996'; E='  -7.95%  '},
    @{Row=12; D='''1.001'; E='  +0.12%  '},
    @{Row=13; D='''5.933'; E='  -6.89%  '},
    @{Row=14; D='''19.24'; E='  -10.35%  '},
    @{Row=15; D='''6.566'; E='  -7.50%  '},
    @{Row=16; D='''1.644.98'; E='  -3.88%  '},
    @{Row=17; D='''0.00001031'; E='  -9.44%  '},
    @{Row=18; D='''0.06542'; E='  -3.24%  '},
    @{Row=19; D='''1.000'; E='  +0.30%  '},
    @{Row=20; D='''77.75'; E='  -8.97%  '},
    @{Row=21; D='''5.908'; E='  -8.19%  '},
    @{Row=22; D='''15.50'; E='  -11.17%  '},
    @{Row=23; D='''12.48'; E='  -6.37%  '},
    @{Row=24; D='''24.371.24'; E='  -2.22%  '},
    @{Row=25; D='''2.454'; E='  -0.17%  '},
    @{Row=26; D='''2.298'; E='  -18.59%  '},
    @{Row=27; D='''145.81'; E='  -3.77%  '},
    @{Row=28; D='''18.49'; E='  -10.43%  '},
    @{Row=29; D='''1.826.10'; E='  -3.96%  '},
    @{Row=30; D='''123.64'; E='  -7.65%  '},
    @{Row=31; D='''1.173'; E='  -6.97%  '},
    @{Row=32; D='''4.066'; E='  -4.02%  '},
    @{Row=33; D='''5.624'; E='  -19.22%  '},
    @{Row=34; D=$null; E='  -5.47%  '},
    @{Row=35; D='''1.641'; E='  -9.67%  '},
    @{Row=36; D='''12.04'; E='  -14.04%  '},
    @{Row=37; D='''1.260'; E='  -2.17%  '},
    @{Row=38; D='''5.152'; E='  -9.30%  '},
    @{Row=39; D='''0.05970'; E='  -11.11%  '},
    @{Row=40; D='''0.02205'; E='  -9.18%  '},
    @{Row=41; D='''0.2041'; E='  -9.46%  '},
    @{Row=42; D='''8.029'; E='  -14.84%  '},
    @{Row=43; D='''1.000'; E='  +0.28%  '},
    @{Row=44; D='''0.5849'; E='  -10.08%  '},
    @{Row=45; D='''3.749'; E='  -2.50%  '},
    @{Row=46; D='''12.44'; E='  -11.76%  '},
    @{Row=47; D='''0.5566'; E='  -10.35%  '},
    @{Row=48; D='''121.74'; E='  -7.00%  '},
    @{Row=49; D='''1.930'; E='  -10.38%  '},
    @{Row=50; D='''0.06889'; E='  -6.08%  '},
    @{Row=51; D='''1.172'; E='  -4.41%  '}
)

foreach ($item in $data) {
    if ($item.D -ne $null) {
        $ws.Cells.Item($item.Row, 4).Value = $item.D
        $ws.Cells.Item($item.Row, 4).Style = "Normal"
    }
    if ($item.E -ne $null) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E
    }
}
